$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Commit: "Renamed control, reduced to single"
#
# The two "control" sample rows (55: CCLB/Vero cells, 56: E6/Vero cells)
# are collapsed into a single row 55. The control is renamed to
# "USA-WA1-2020-TCE", the exp id becomes "VSP0002", and the surviving
# numeric/metadata columns are the ones that used to live on row 56
# (composite, NP-OP, 30.02 / 99.8 / 99.8).

# Text columns whose literal content happens to look like a date
# ("2020-03-28") or a pure digit string ("20200328") get auto-converted
# by Excel's normal Value-assignment type-sniffer (to a date serial, or
# a quote-prefixed cell carrying a spurious style). To land them as
# plain shared-string text - exactly like every other sampleDate /
# sampleDate2 cell in this sheet - stage the literal text as a formula
# result in a scratch cell and transplant it with a values-only paste,
# which copies the already-resolved text without re-running type
# inference.
$ws.Range("Z1").Formula = '=TEXT("2020-03-28","@")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("F55").PasteSpecial(-4163)

$ws.Range("Z1").Formula = '=TEXT("20200328","@")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("J55").PasteSpecial(-4163)

$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = $false

# Plain labels/numbers - safe to set directly, no type-sniffing surprises.
$ws.Range("A55").Value = "USA-WA1-2020-TCE"
$ws.Range("B55").Value = "VSP0002"
$ws.Range("C55").Value = "composite"
$ws.Range("E55").Value = "NP-OP"
$ws.Range("G55").Value = 30.02
$ws.Range("H55").Value = 99.8
$ws.Range("I55").Value = 99.8
$ws.Range("K55").Value = "USA-WA1-2020-TCE NP-OP 20200328"

# Row 56 (the other old control row) is gone now that the two rows were
# reduced to one.
$ws.Rows.Item(56).Delete()
